$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 11) down across
# the new rows (12-24) so the new cells pick up the same cell style (s="1")
# as all the other data rows.
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C24").PasteSpecial(-4122)
$ws.Rows("12:24").RowHeight = 15.75

# --- Column C (Syntax_Template) first, so the new shared strings for the
# syntax templates are appended to sharedStrings.xml before the new
# keyword/category labels (matches original authoring order). ---
$ws.Range("C12").Value = "DESCRIPTIVES VARIABLES=var /STATISTICS=MEAN STDDEV MIN MAX."
$ws.Range("C13").Value = "DESCRIPTIVES VARIABLES=var /SAVE."
$ws.Range("C14").Value = "T-TEST GROUPS=group(1 2) /VARIABLES=test_var."
$ws.Range("C15").Value = "T-TEST PAIRS=Before WITH After (PAIRED)."
$ws.Range("C16").Value = "ONEWAY var BY group /STATISTICS DESCRIPTIVES /POSTHOC=TUKEY."
$ws.Range("C17").Value = "NPAR TESTS /K-S(NORMAL)=var."
$ws.Range("C18").Value = "CORRELATIONS /VARIABLES=var1 var2 /PRINT=TWOTAIL."
$ws.Range("C19").Value = "REGRESSION /DEPENDENT=y /METHOD=ENTER x."
$ws.Range("C20").Value = "REGRESSION /DEPENDENT=y /METHOD=ENTER x1 x2 x3."
$ws.Range("C21").Value = "GRAPH /BAR(MEAN)=var BY group."
$ws.Range("C22").Value = "GRAPH /PIE=SUM(var) BY group."
$ws.Range("C23").Value = "GRAPH /HISTOGRAM(NORMAL)=var."
$ws.Range("C24").Value = "RECODE var (Low thru 1990=1) (1991 thru Hi=2) INTO newvar."

# --- Columns A (Keyword) and B (Category) share the same label text for
# each new row. ---
$ws.Range("A12").Value = "Descriptives"
$ws.Range("B12").Value = "Descriptives"
$ws.Range("A13").Value = "Z-Scores"
$ws.Range("B13").Value = "Z-Scores"
$ws.Range("A14").Value = "Independent T-Test"
$ws.Range("B14").Value = "Independent T-Test"
$ws.Range("A15").Value = "Paired T-Test"
$ws.Range("B15").Value = "Paired T-Test"
$ws.Range("A16").Value = "One-Way ANOVA"
$ws.Range("B16").Value = "One-Way ANOVA"
$ws.Range("A17").Value = "Normality (K-S)"
$ws.Range("B17").Value = "Normality (K-S)"
$ws.Range("A18").Value = "Correlation"
$ws.Range("B18").Value = "Correlation"
$ws.Range("A19").Value = "Simple Regression"
$ws.Range("B19").Value = "Simple Regression"
$ws.Range("A20").Value = "Multiple Regression"
$ws.Range("B20").Value = "Multiple Regression"
$ws.Range("A21").Value = "Bar Chart"
$ws.Range("B21").Value = "Bar Chart"
$ws.Range("A22").Value = "Pie Chart"
$ws.Range("B22").Value = "Pie Chart"
$ws.Range("A23").Value = "Histogram"
$ws.Range("B23").Value = "Histogram"
$ws.Range("A24").Value = "Recode Data"
$ws.Range("B24").Value = "Recode Data"

# Match the final view state from the authored file: sheet active/selected,
# with the newly added rows selected.
$ws.Select()
[void]$ws.Range("A12:A24").Select()
